$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Values are set with a leading quote-prefix via .Formula so that
# numeric-looking strings (e.g. "540.81") stay stored as text (inlineStr/shared
# string), matching the original cell type. The quote-prefix style applied by
# Excel is then reset via .Style so no stray style index is introduced.

$cell = $ws.Cells.Item(2, 4)
$cell.Formula = "'57.924.34"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.Formula = "'  +0.48%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.Formula = "'2.346.63"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.Formula = "'  +0.88%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 5)
$cell.Formula = "'  -0.15%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.Formula = "'540.81"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.Formula = "'  -0.12%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 4)
$cell.Formula = "'134.50"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.Formula = "'  -0.13%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 5)
$cell.Formula = "'  +0.31%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(8, 4)
$cell.Formula = "'0.571"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.Formula = "'  +6.78%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 5)
$cell.Formula = "'  +0.47%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.Formula = "'5.52"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.Formula = "'  +2.83%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 5)
$cell.Formula = "'  -1.63%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 5)
$cell.Formula = "'  +1.21%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.Formula = "'2.765.15"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.Formula = "'  -0.51%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.Formula = "'57.857.78"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.Formula = "'  +0.24%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 5)
$cell.Formula = "'  +0.86%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 4)
$cell.Formula = "'2.351.51"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.Formula = "'  +0.32%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.Formula = "'10.69"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.Formula = "'  +1.53%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 5)
$cell.Formula = "'  +2.31%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.Formula = "'329.60"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.Formula = "'  -2.43%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 4)
$cell.Formula = "'6.70"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.Formula = "'  -0.83%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.Formula = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.Formula = "'  +0.08%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.Formula = "'62.85"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.Formula = "'  +1.36%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.Formula = "'0.164"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.Formula = "'  -2.95%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 4)
$cell.Formula = "'0.997"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.Formula = "'  +0.26%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 4)
$cell.Formula = "'8.34"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.Formula = "'  -1.41%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 4)
$cell.Formula = "'1.34"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.Formula = "'  -5.74%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 5)
$cell.Formula = "'  +0.42%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 4)
$cell.Formula = "'169.99"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.Formula = "'  -0.22%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 5)
$cell.Formula = "'  +0.14%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 5)
$cell.Formula = "'  -0.50%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 4)
$cell.Formula = "'1.02"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.Formula = "'  +0.23%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 5)
$cell.Formula = "'  -0.94%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 5)
$cell.Formula = "'  -0.03%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 4)
$cell.Formula = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.Formula = "'  +0.58%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 4)
$cell.Formula = "'4.19"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.Formula = "'  +1.81%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 4)
$cell.Formula = "'1.23"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.Formula = "'  -1.67%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 5)
$cell.Formula = "'  +0.11%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 4)
$cell.Formula = "'39.10"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.Formula = "'  -0.60%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 4)
$cell.Formula = "'142.25"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.Formula = "'  -4.04%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 4)
$cell.Formula = "'0.377"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.Formula = "'  -0.45%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 5)
$cell.Formula = "'  +0.78%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 4)
$cell.Formula = "'288.77"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.Formula = "'  +2.15%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 4)
$cell.Formula = "'0.0948"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.Formula = "'  +1.88%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 5)
$cell.Formula = "'  +0.83%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 4)
$cell.Formula = "'19.14"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.Formula = "'  -0.15%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 5)
$cell.Formula = "'  +1.36%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 5)
$cell.Formula = "'  +1.48%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 5)
$cell.Formula = "'  -0.37%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 5)
$cell.Formula = "'  +0.52%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 2)
$cell.Formula = "'BitgetToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.Formula = "'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.Formula = "'0.953"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.Formula = "'  +0.89%  "
$cell.Style = "Normal"
